# Navigate to cart started... working on filling out address information

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Address")

# ---- Duplicate the existing header/data row pair (rows 1-2) down to rows 4-5,
#      inserting-with-copy so the new cells inherit the existing cell style (style 1). ----
$ws.Range("A1:F2").Copy()
$ws.Range("A4:F5").Insert(-4121)

# ---- Insert two new columns (G, H) across both blocks (rows 1-5) in one shot so the
#      new cells inherit the existing column formatting (style 1) too. ----
$ws.Range("F1:F5").Copy()
$ws.Range("G1:G5").Insert(-4161)
$ws.Range("F1:F5").Copy()
$ws.Range("H1:H5").Insert(-4161)

# H2/H5 ("Company Name" data cells) should end up unstyled (no explicit style).
$ws.Cells.Item(2,8).ClearFormats()
$ws.Cells.Item(5,8).ClearFormats()

# ---- Row 1 : Billing Address header block ----
$ws.Cells.Item(1,1).Value = "Billing Address"
$ws.Cells.Item(1,2).Value = "First Name"
$ws.Cells.Item(1,3).Value = "Last Name"
$ws.Cells.Item(1,4).Value = "Street Address"
$ws.Cells.Item(1,5).Value = "ZIP Code"
$ws.Cells.Item(1,6).Value = "Area Code"
$ws.Cells.Item(1,7).Value = "Primary Phone"
$ws.Cells.Item(1,8).Value = "Company Name"

# ---- Row 2 : Billing Address data ----
$ws.Cells.Item(2,1).Value = ""
$ws.Cells.Item(2,2).Value = "Steve"
$ws.Cells.Item(2,3).Value = "Jobs"
$ws.Cells.Item(2,4).Value = "1111 First Street"
$ws.Cells.Item(2,5).Value = 76013
$ws.Cells.Item(2,6).Value = 208
$ws.Cells.Item(2,7).Value = 5554970
$ws.Cells.Item(2,8).Value = "Apple"

# ---- Row 4 : Shipping Address header block ----
$ws.Cells.Item(4,1).Value = "Shipping Address"
$ws.Cells.Item(4,2).Value = "First Name"
$ws.Cells.Item(4,3).Value = "Last Name"
$ws.Cells.Item(4,4).Value = "Street Address"
$ws.Cells.Item(4,5).Value = "ZIP Code"
$ws.Cells.Item(4,6).Value = "Area Code"
$ws.Cells.Item(4,7).Value = "Primary Phone"
$ws.Cells.Item(4,8).Value = "Company Name"

# ---- Row 5 : Shipping Address data ----
$ws.Cells.Item(5,1).Value = ""
$ws.Cells.Item(5,2).Value = "Steve"
$ws.Cells.Item(5,3).Value = "Jobs"
$ws.Cells.Item(5,4).Value = "1111 First Street"
$ws.Cells.Item(5,5).Value = 76013
$ws.Cells.Item(5,6).Value = 208
$ws.Cells.Item(5,7).Value = 5554970
$ws.Cells.Item(5,8).Value = "Apple"

# ---- Merge the section-title cell with the blank cell below it (do this BEFORE
#      applying distinct per-cell formatting, since Merge() normalizes the whole
#      merged range to a single style). ----
$ws.Range("A1:A2").Merge()
$ws.Range("A4:A5").Merge()

# ---- Formatting: bold + centered "section" header cells ----
$ws.Range("A1").Font.Bold = $true
$ws.Range("A1").HorizontalAlignment = -4108
$ws.Range("A2").HorizontalAlignment = -4108
$ws.Range("A4").Font.Bold = $true
$ws.Range("A4").HorizontalAlignment = -4108
$ws.Range("A5").Font.Bold = $true
$ws.Range("A5").HorizontalAlignment = -4108

# ---- Page setup: portrait orientation (creates printer-settings relationship) ----
$ws.PageSetup.Orientation = 1

# ---- Selection / active tab: the Address sheet becomes the active / visible sheet ----
$ws.Range("F7").Select()
$ws.Activate()

$wb.Save()
